$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count
    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            $val = $cell.Value()
            if ($val -ne $null -and $val -is [string]) {
                $newVal = $val -replace '\[\d+\]', ''
                $newVal = $newVal -replace "`n", ' '
                if ($newVal -ne $val) {
                    $cell.Value = $newVal
                }
            }
        }
    }
}
